$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-6 (years 2014-2018) with corrected values

# Row 2
$ws.Range("D2").Value = 3270
$ws.Range("E2").Value = 251
$ws.Range("F2").Value = 251
$ws.Range("G2").Value = 243
$ws.Range("H2").Value = 164
$ws.Range("I2").Value = 166
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 4018
$ws.Range("L2").Value = 1277
$ws.Range("M2").Value = 2741
$ws.Range("N2").Value = 2553
$ws.Range("O2").Value = 188
$ws.Range("P2").Value = 183
$ws.Range("Q2").Value = 419
$ws.Range("R2").Value = -337
$ws.Range("S2").Value = -82
$ws.Range("T2").Value = 69
$ws.Range("U2").Value = 350
$ws.Range("V2").Value = 701
$ws.Range("W2").Value = 7.68
$ws.Range("X2").Value = 5.02
$ws.Range("Y2").Value = 6.68
$ws.Range("Z2").Value = 4.13
$ws.Range("AA2").Value = 46.57
$ws.Range("AB2").Value = 1313.88
$ws.Range("AC2").Value = 4588
$ws.Range("AD2").Value = 7.35
$ws.Range("AE2").Value = 69649
$ws.Range("AF2").Value = 0.48
$ws.Range("AG2").Value = 375
$ws.Range("AH2").Value = 1.11
$ws.Range("AI2").Value = 8.300000000000001
$ws.Range("AJ2").Value = 3665813

# Row 3
$ws.Range("D3").Value = 2393
$ws.Range("E3").Value = 83
$ws.Range("F3").Value = 87
$ws.Range("G3").Value = -299
$ws.Range("H3").Value = -235
$ws.Range("I3").Value = -214
$ws.Range("J3").Value = -20
$ws.Range("K3").Value = 3756
$ws.Range("L3").Value = 1216
$ws.Range("M3").Value = 2540
$ws.Range("N3").Value = 2373
$ws.Range("O3").Value = 167
$ws.Range("P3").Value = 193
$ws.Range("Q3").Value = 206
$ws.Range("R3").Value = 82
$ws.Range("S3").Value = -284
$ws.Range("T3").Value = 31
$ws.Range("U3").Value = 175
$ws.Range("V3").Value = 396
$ws.Range("W3").Value = 3.47
$ws.Range("X3").Value = -9.800000000000001
$ws.Range("Y3").Value = -8.699999999999999
$ws.Range("Z3").Value = -6.03
$ws.Range("AA3").Value = 47.88
$ws.Range("AB3").Value = 1146.28
$ws.Range("AC3").Value = -5671
$ws.Range("AD3").Value = -4.48
$ws.Range("AE3").Value = 61323
$ws.Range("AF3").Value = 0.41
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 3869237

# Row 4
$ws.Range("D4").Value = 2525
$ws.Range("E4").Value = 82
$ws.Range("F4").Value = 82
$ws.Range("G4").Value = 297
$ws.Range("H4").Value = 170
$ws.Range("I4").Value = 156
$ws.Range("J4").Value = 14
$ws.Range("K4").Value = 3627
$ws.Range("L4").Value = 978
$ws.Range("M4").Value = 2649
$ws.Range("N4").Value = 2519
$ws.Range("O4").Value = 131
$ws.Range("P4").Value = 193
$ws.Range("Q4").Value = -19
$ws.Range("R4").Value = -107
$ws.Range("S4").Value = 162
$ws.Range("T4").Value = 324
$ws.Range("U4").Value = -343
$ws.Range("V4").Value = 277
$ws.Range("W4").Value = 3.26
$ws.Range("X4").Value = 6.73
$ws.Range("Y4").Value = 6.37
$ws.Range("Z4").Value = 4.6
$ws.Range("AA4").Value = 36.92
$ws.Range("AB4").Value = 1225.03
$ws.Range("AC4").Value = 4028
$ws.Range("AD4").Value = 6.03
$ws.Range("AE4").Value = 65091
$ws.Range("AF4").Value = 0.37
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 3869237

# Row 5
$ws.Range("D5").Value = 5128
$ws.Range("E5").Value = 288
$ws.Range("F5").Value = 288
$ws.Range("G5").Value = 352
$ws.Range("H5").Value = 272
$ws.Range("I5").Value = 273
$ws.Range("J5").Value = -2
$ws.Range("K5").Value = 7258
$ws.Range("L5").Value = 3102
$ws.Range("M5").Value = 4156
$ws.Range("N5").Value = 2753
$ws.Range("O5").Value = 1403
$ws.Range("P5").Value = 201
$ws.Range("Q5").Value = 191
$ws.Range("R5").Value = -598
$ws.Range("S5").Value = 418
$ws.Range("T5").Value = 240
$ws.Range("U5").Value = -48
$ws.Range("V5").Value = 1836
$ws.Range("W5").Value = 5.61
$ws.Range("X5").Value = 5.3
$ws.Range("Y5").Value = 10.37
$ws.Range("Z5").Value = 4.99
$ws.Range("AA5").Value = 74.65000000000001
$ws.Range("AB5").Value = 1265.8
$ws.Range("AC5").Value = 7032
$ws.Range("AD5").Value = 4.96
$ws.Range("AE5").Value = 68314
$ws.Range("AF5").Value = 0.51
$ws.Range("AG5").Value = 375
$ws.Range("AH5").Value = 1.08
$ws.Range("AI5").Value = 5.53
$ws.Range("AJ5").Value = 4029782

# Row 6
$ws.Range("D6").Value = 7094
$ws.Range("E6").Value = 1166
$ws.Range("F6").Value = 1166
$ws.Range("G6").Value = 1215
$ws.Range("H6").Value = 984
$ws.Range("I6").Value = 793
$ws.Range("K6").Value = 7468
$ws.Range("L6").Value = 2403
$ws.Range("M6").Value = 5065
$ws.Range("N6").Value = 3644
$ws.Range("P6").Value = 201
$ws.Range("Q6").Value = 1373
$ws.Range("R6").Value = -657
$ws.Range("S6").Value = -737
$ws.Range("T6").Value = 349
$ws.Range("U6").Value = 1024
$ws.Range("V6").Value = 1208
$ws.Range("W6").Value = 16.44
$ws.Range("X6").Value = 13.88
$ws.Range("Y6").Value = 24.79
$ws.Range("Z6").Value = 13.37
$ws.Range("AA6").Value = 47.44
$ws.Range("AB6").Value = 1708.28
$ws.Range("AC6").Value = 19674
$ws.Range("AD6").Value = 3.01
$ws.Range("AE6").Value = 90415
$ws.Range("AF6").Value = 0.66
$ws.Range("AG6").Value = 750
$ws.Range("AH6").Value = 1.26
$ws.Range("AI6").Value = 3.81
$ws.Range("AJ6").Value = 4029782

# Rows 7-9 (2019E/2020E/2021E) lose all data cells (D:AJ), keeping only A, B, C
$ws.Range("D7:AJ9").ClearContents()
